$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 14 entirely - shifts all rows below it up by one.
$ws.Rows.Item(14).Delete()

# Restore the view state (active cell) to match where the author ended up
# after the edit.
[void]$ws.Range("A14").Select()
